$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Cells.Item(1,4).Value = "Jan_2026"
$ws.Cells.Item(1,5).Value = "Dec_2025"
$ws.Cells.Item(1,6).Value = "Nov_2025"

# Update data rows 2-19 with new values (row order/content reshuffled)
$ws.Cells.Item(2,1).Value = "INE018A01030"
$ws.Cells.Item(2,2).Value = "Larsen & Toubro Limited"
$ws.Cells.Item(2,4).Value = 9.769026
$ws.Cells.Item(2,5).Value = 9.405836
$ws.Cells.Item(2,6).Value = 7.138231
$ws.Cells.Item(2,7).Value = 0.3631899999999995
$ws.Cells.Item(2,8).Value = 2.630795

$ws.Cells.Item(3,1).Value = "INE040A01034"
$ws.Cells.Item(3,2).Value = "HDFC Bank Limited"
$ws.Cells.Item(3,4).Value = 9.010578
$ws.Cells.Item(3,5).Value = 0.878126
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 8.132452
$ws.Cells.Item(3,8).Value = 9.010578

$ws.Cells.Item(4,1).Value = "INE423A01024"
$ws.Cells.Item(4,2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(4,4).Value = 8.730519
$ws.Cells.Item(4,5).Value = 8.973312
$ws.Cells.Item(4,6).Value = 8.706628
$ws.Cells.Item(4,7).Value = -0.2427930000000007
$ws.Cells.Item(4,8).Value = 0.023890999999999

$ws.Cells.Item(5,1).Value = "INE364U01010"
$ws.Cells.Item(5,2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(5,4).Value = 8.249091
$ws.Cells.Item(5,5).Value = 9.109251
$ws.Cells.Item(5,6).Value = 7.106706
$ws.Cells.Item(5,7).Value = -0.8601600000000005
$ws.Cells.Item(5,8).Value = 1.142385

$ws.Cells.Item(6,1).Value = "INE180C01042"
$ws.Cells.Item(6,2).Value = "Capri Global Capital Limited"
$ws.Cells.Item(6,4).Value = 7.92449
$ws.Cells.Item(6,5).Value = 7.644202
$ws.Cells.Item(6,6).Value = 7.617549
$ws.Cells.Item(6,7).Value = 0.2802879999999996
$ws.Cells.Item(6,8).Value = 0.3069409999999992

$ws.Cells.Item(7,1).Value = "INE090A01021"
$ws.Cells.Item(7,2).Value = "ICICI Bank Limited"
$ws.Cells.Item(7,4).Value = 6.625112
$ws.Cells.Item(7,5).Value = 6.087764
$ws.Cells.Item(7,6).Value = 6.097153
$ws.Cells.Item(7,7).Value = 0.5373479999999997
$ws.Cells.Item(7,8).Value = 0.5279590000000001

$ws.Cells.Item(8,1).Value = "INE775A01035"
$ws.Cells.Item(8,2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(8,4).Value = 6.249481
$ws.Cells.Item(8,5).Value = 6.155655
$ws.Cells.Item(8,6).Value = 0
$ws.Cells.Item(8,7).Value = 0.09382599999999996
$ws.Cells.Item(8,8).Value = 6.249481

$ws.Cells.Item(9,1).Value = "INE795G01014"
$ws.Cells.Item(9,2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(9,4).Value = 6.150729
$ws.Cells.Item(9,5).Value = 5.849838
$ws.Cells.Item(9,6).Value = 0
$ws.Cells.Item(9,7).Value = 0.300891
$ws.Cells.Item(9,8).Value = 6.150729

$ws.Cells.Item(10,1).Value = "INE237A01036"
$ws.Cells.Item(10,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(10,4).Value = 5.930867
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 0
$ws.Cells.Item(10,7).Value = 5.930867
$ws.Cells.Item(10,8).Value = 5.930867

$ws.Cells.Item(11,1).Value = "INE917I01010"
$ws.Cells.Item(11,2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(11,4).Value = 4.703453
$ws.Cells.Item(11,5).Value = 4.245269
$ws.Cells.Item(11,6).Value = 3.992703
$ws.Cells.Item(11,7).Value = 0.4581839999999993
$ws.Cells.Item(11,8).Value = 0.7107499999999995

$ws.Cells.Item(12,1).Value = "INE237A01028"
$ws.Cells.Item(12,2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).Value = 5.933186
$ws.Cells.Item(12,6).Value = 5.545717
$ws.Cells.Item(12,7).Value = -5.933186
$ws.Cells.Item(12,8).Value = -5.545717

$ws.Cells.Item(13,1).Value = "INE245A01021"
$ws.Cells.Item(13,2).Value = "Tata Power Company Limited"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 8.23647
$ws.Cells.Item(13,7).Value = 0
$ws.Cells.Item(13,8).Value = -8.23647

$ws.Cells.Item(14,1).Value = "INE271C01023"
$ws.Cells.Item(14,2).Value = "DLF Limited"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 6.855223
$ws.Cells.Item(14,6).Value = 6.988498
$ws.Cells.Item(14,7).Value = -6.855223
$ws.Cells.Item(14,8).Value = -6.988498

$ws.Cells.Item(15,1).Value = "INE044A01036"
$ws.Cells.Item(15,2).Value = "Sun Pharmaceutical Industries Limited"
$ws.Cells.Item(15,4).Value = 0
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 3.59477
$ws.Cells.Item(15,7).Value = 0
$ws.Cells.Item(15,8).Value = -3.59477

$ws.Cells.Item(16,1).Value = "INE423A20016"
$ws.Cells.Item(16,2).Value = "Adani Enterprises Limited Rights"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 0.202705
$ws.Cells.Item(16,7).Value = 0
$ws.Cells.Item(16,8).Value = -0.202705

$ws.Cells.Item(17,1).Value = "INE669C01036"
$ws.Cells.Item(17,2).Value = "Tech Mahindra Limited"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).Value = 4.001813
$ws.Cells.Item(17,6).Value = 2.92597
$ws.Cells.Item(17,7).Value = -4.001813
$ws.Cells.Item(17,8).Value = -2.92597

$ws.Cells.Item(18,1).Value = "INE726G01019"
$ws.Cells.Item(18,2).Value = "ICICI Prudential Life Insurance Co Ltd"
$ws.Cells.Item(18,4).Value = 0
$ws.Cells.Item(18,5).Value = 3.030567
$ws.Cells.Item(18,6).Value = 0
$ws.Cells.Item(18,7).Value = -3.030567
$ws.Cells.Item(18,8).Value = 0

$ws.Cells.Item(19,1).Value = "INE918I01026"
$ws.Cells.Item(19,2).Value = "Bajaj Finserv Ltd."
$ws.Cells.Item(19,4).Value = 0
$ws.Cells.Item(19,5).Value = 1.355082
$ws.Cells.Item(19,6).Value = 1.347121
$ws.Cells.Item(19,7).Value = -1.355082
$ws.Cells.Item(19,8).Value = -1.347121

# Remove now-obsolete rows 20-22
$ws.Range("A20:H22").Delete()
